$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 5701.5557
$ws.Cells.Item(6, 9).Value = 10090
$ws.Cells.Item(6, 11).Value = 30270
$ws.Cells.Item(6, 13).Value = -30158
$ws.Cells.Item(33, 8).Value = 151.58333
$ws.Cells.Item(33, 9).Value = 86.09999999999999
$ws.Cells.Item(33, 10).Value = 479
$ws.Cells.Item(33, 11).Value = 86.09999999999999
$ws.Cells.Item(33, 12).Value = 479
$ws.Cells.Item(33, 13).Value = 142.9
$ws.Cells.Item(33, 14).Value = -937
$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 13).Value = ""
$ws.Cells.Item(98, 8).Value = 2862.7307
$ws.Cells.Item(98, 9).Value = 2862.7307
$ws.Cells.Item(98, 11).Value = 2862.7307
$ws.Cells.Item(98, 13).Value = -1364.7307
$ws.Cells.Item(122, 8).Value = 2862.7307
$ws.Cells.Item(122, 9).Value = 2862.7307
$ws.Cells.Item(122, 11).Value = 8588.1921
$ws.Cells.Item(122, 13).Value = -6138.1921
$ws.Cells.Item(135, 8).Value = 1122.375
$ws.Cells.Item(135, 9).Value = 292.33334
$ws.Cells.Item(135, 11).Value = 2631.00006
$ws.Cells.Item(135, 13).Value = -96.0000600000003
$ws.Cells.Item(138, 8).Value = 1591.8036
$ws.Cells.Item(138, 9).Value = 1405.8462
$ws.Cells.Item(138, 10).Value = 1752.9667
$ws.Cells.Item(138, 11).Value = 4217.5386
$ws.Cells.Item(138, 12).Value = 5258.9001
$ws.Cells.Item(138, 13).Value = 922.4614000000001
$ws.Cells.Item(138, 14).Value = -15538.9001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 745.1818
$ws.Cells.Item(61, 9).Value = 721.8889
$ws.Cells.Item(61, 10).Value = 850
$ws.Cells.Item(61, 11).Value = 721.8889
$ws.Cells.Item(61, 12).Value = 850
$ws.Cells.Item(61, 13).Value = -509.8889
$ws.Cells.Item(61, 14).Value = -1274
$ws.Cells.Item(132, 8).Value = 1864.2333
$ws.Cells.Item(132, 9).Value = 1572.619
$ws.Cells.Item(132, 10).Value = 2544.6667
$ws.Cells.Item(132, 11).Value = 4717.857
$ws.Cells.Item(132, 12).Value = 7634.000100000001
$ws.Cells.Item(132, 13).Value = -2187.857
$ws.Cells.Item(132, 14).Value = -12694.0001
$ws.Cells.Item(136, 8).Value = 745.1818
$ws.Cells.Item(136, 9).Value = 721.8889
$ws.Cells.Item(136, 10).Value = 850
$ws.Cells.Item(136, 11).Value = 2165.6667
$ws.Cells.Item(136, 12).Value = 2550
$ws.Cells.Item(136, 13).Value = 384.3332999999998
$ws.Cells.Item(136, 14).Value = -7650
$ws.Cells.Item(141, 8).Value = 30429
$ws.Cells.Item(141, 10).Value = 30429
$ws.Cells.Item(141, 12).Value = 30429
$ws.Cells.Item(141, 14).Value = -40789

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(112, 8).Value = 100000000
$ws.Cells.Item(112, 9).Value = 100000000
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 11).Value = 100000000
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 13).Value = -99998523
$ws.Cells.Item(112, 14).Value = ""
$ws.Cells.Item(132, 8).Value = 54999
$ws.Cells.Item(132, 10).Value = 54999
$ws.Cells.Item(132, 12).Value = 54999
$ws.Cells.Item(132, 14).Value = -65119
$ws.Cells.Item(134, 8).Value = 6210.5
$ws.Cells.Item(134, 9).Value = 952.65
$ws.Cells.Item(134, 11).Value = 2857.95
$ws.Cells.Item(134, 13).Value = -322.9499999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1344.1428
$ws.Cells.Item(31, 9).Value = 1094.619
$ws.Cells.Item(31, 10).Value = 2092.7144
$ws.Cells.Item(31, 11).Value = 1094.619
$ws.Cells.Item(31, 12).Value = 2092.7144
$ws.Cells.Item(31, 13).Value = -799.6189999999999
$ws.Cells.Item(31, 14).Value = -2682.7144
$ws.Cells.Item(34, 8).Value = 1344.1428
$ws.Cells.Item(34, 9).Value = 1094.619
$ws.Cells.Item(34, 10).Value = 2092.7144
$ws.Cells.Item(34, 11).Value = 1094.619
$ws.Cells.Item(34, 12).Value = 2092.7144
$ws.Cells.Item(34, 13).Value = -892.6189999999999
$ws.Cells.Item(34, 14).Value = -2496.7144
$ws.Cells.Item(58, 8).Value = 2248.6155
$ws.Cells.Item(58, 9).Value = 1803.5555
$ws.Cells.Item(58, 10).Value = 3250
$ws.Cells.Item(58, 11).Value = 1803.5555
$ws.Cells.Item(58, 12).Value = 3250
$ws.Cells.Item(58, 13).Value = -1600.5555
$ws.Cells.Item(58, 14).Value = -3656
$ws.Cells.Item(136, 8).Value = 2248.6155
$ws.Cells.Item(136, 9).Value = 1803.5555
$ws.Cells.Item(136, 10).Value = 3250
$ws.Cells.Item(136, 11).Value = 5410.666499999999
$ws.Cells.Item(136, 12).Value = 9750
$ws.Cells.Item(136, 13).Value = -2860.666499999999
$ws.Cells.Item(136, 14).Value = -14850

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 604.65515
$ws.Cells.Item(5, 9).Value = 603.5217
$ws.Cells.Item(5, 10).Value = 609
$ws.Cells.Item(5, 11).Value = 1810.5651
$ws.Cells.Item(5, 12).Value = 1827
$ws.Cells.Item(5, 13).Value = -1698.5651
$ws.Cells.Item(5, 14).Value = -2051
$ws.Cells.Item(7, 8).Value = 383.92307
$ws.Cells.Item(7, 9).Value = 368.9
$ws.Cells.Item(7, 10).Value = 434
$ws.Cells.Item(7, 11).Value = 1106.7
$ws.Cells.Item(7, 12).Value = 1302
$ws.Cells.Item(7, 13).Value = -994.6999999999998
$ws.Cells.Item(7, 14).Value = -1526
$ws.Cells.Item(80, 8).Value = 4570
$ws.Cells.Item(80, 10).Value = 4570
$ws.Cells.Item(80, 12).Value = 13710
$ws.Cells.Item(80, 14).Value = -15582
$ws.Cells.Item(83, 8).Value = 4570
$ws.Cells.Item(83, 10).Value = 4570
$ws.Cells.Item(83, 12).Value = 41130
$ws.Cells.Item(83, 14).Value = -50490
$ws.Cells.Item(92, 8).Value = 828
$ws.Cells.Item(92, 10).Value = 800
$ws.Cells.Item(92, 12).Value = 2400
$ws.Cells.Item(92, 14).Value = -4896
$ws.Cells.Item(122, 8).Value = 741.7778
$ws.Cells.Item(122, 9).Value = 517.2
$ws.Cells.Item(122, 11).Value = 4654.8
$ws.Cells.Item(122, 13).Value = -2204.8
$ws.Cells.Item(131, 8).Value = 12346853
$ws.Cells.Item(131, 9).Value = 333333600
$ws.Cells.Item(131, 10).Value = 1208.6794
$ws.Cells.Item(131, 11).Value = 1000000800
$ws.Cells.Item(131, 12).Value = 3626.0382
$ws.Cells.Item(131, 13).Value = -999995760
$ws.Cells.Item(131, 14).Value = -13706.0382
$ws.Cells.Item(135, 8).Value = 604.65515
$ws.Cells.Item(135, 9).Value = 603.5217
$ws.Cells.Item(135, 10).Value = 609
$ws.Cells.Item(135, 11).Value = 5431.6953
$ws.Cells.Item(135, 12).Value = 5481
$ws.Cells.Item(135, 13).Value = -2896.6953
$ws.Cells.Item(135, 14).Value = -10551

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 5115600
$ws.Cells.Item(12, 9).Value = 5194647
$ws.Cells.Item(12, 10).Value = 4667665
$ws.Cells.Item(12, 11).Value = 5194647
$ws.Cells.Item(12, 12).Value = 4667665
$ws.Cells.Item(12, 13).Value = -5194507
$ws.Cells.Item(12, 14).Value = -4667945
$ws.Cells.Item(70, 8).Value = 150001340
$ws.Cells.Item(70, 10).Value = 100002010
$ws.Cells.Item(70, 12).Value = 100002010
$ws.Cells.Item(70, 14).Value = -100002550
$ws.Cells.Item(73, 8).Value = 150001340
$ws.Cells.Item(73, 10).Value = 100002010
$ws.Cells.Item(73, 12).Value = 100002010
$ws.Cells.Item(73, 14).Value = -100003882
$ws.Cells.Item(138, 8).Value = 34857.25
$ws.Cells.Item(138, 10).Value = 34857.25
$ws.Cells.Item(138, 12).Value = 34857.25
$ws.Cells.Item(138, 14).Value = -45137.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 21505000
$ws.Cells.Item(5, 10).Value = 21505000
$ws.Cells.Item(5, 12).Value = 21505000
$ws.Cells.Item(5, 14).Value = -21505224
$ws.Cells.Item(126, 8).Value = 1462.75
$ws.Cells.Item(126, 9).Value = 1417
$ws.Cells.Item(126, 11).Value = 4251
$ws.Cells.Item(126, 13).Value = -1781
$ws.Cells.Item(127, 8).Value = 67000
$ws.Cells.Item(127, 9).Value = 0
$ws.Cells.Item(127, 11).Value = 0
$ws.Cells.Item(127, 13).Value = ""

